# Scheduled market-data refresh: update computed price/profit cells
# on each job sheet to match the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1238.125
$ws.Range("J17").Value = 1267.3182
$ws.Range("L17").Value = 3801.9546
$ws.Range("N17").Value = -4137.9546
$ws.Range("H18").Value = 3886.1538
$ws.Range("I18").Value = 4892.2
$ws.Range("K18").Value = 4892.2
$ws.Range("M18").Value = -4608.2
$ws.Range("H100").Value = 31437166
$ws.Range("I100").Value = 62750830
$ws.Range("K100").Value = 62750830
$ws.Range("M100").Value = -62750289
$ws.Range("H112").Value = 2334.2666
$ws.Range("I112").Value = 1349.5
$ws.Range("K112").Value = 4048.5
$ws.Range("M112").Value = -2940.5
$ws.Range("H129").Value = 1321
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7919.9585
$ws.Range("I61").Value = 8565.723
$ws.Range("J61").Value = 5982.6665
$ws.Range("K61").Value = 8565.723
$ws.Range("L61").Value = 5982.6665
$ws.Range("M61").Value = -8353.723
$ws.Range("N61").Value = -6406.6665
$ws.Range("H110").Value = 2531.3845
$ws.Range("I110").Value = 1600.8889
$ws.Range("K110").Value = 1600.8889
$ws.Range("M110").Value = 444.1111000000001
$ws.Range("H122").Value = 1158903.6
$ws.Range("I122").Value = 4508.6665
$ws.Range("K122").Value = 13525.9995
$ws.Range("M122").Value = -11075.9995
$ws.Range("H136").Value = 7919.9585
$ws.Range("I136").Value = 8565.723
$ws.Range("J136").Value = 5982.6665
$ws.Range("K136").Value = 25697.169
$ws.Range("L136").Value = 17947.9995
$ws.Range("M136").Value = -23147.169
$ws.Range("N136").Value = -23047.9995

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4242.6875
$ws.Range("I20").Value = 2203.7
$ws.Range("K20").Value = 2203.7
$ws.Range("M20").Value = -1956.7
$ws.Range("H86").Value = 5464.476
$ws.Range("J86").Value = 1486.8182
$ws.Range("L86").Value = 1486.8182
$ws.Range("N86").Value = -3732.8182
$ws.Range("H89").Value = 5464.476
$ws.Range("J89").Value = 1486.8182
$ws.Range("L89").Value = 7434.090999999999
$ws.Range("N89").Value = -18666.091
$ws.Range("H105").Value = 208400
$ws.Range("J105").Value = 10500
$ws.Range("L105").Value = 10500
$ws.Range("N105").Value = -13994
$ws.Range("H107").Value = 2972.8572
$ws.Range("I107").Value = 3301.6667
$ws.Range("K107").Value = 3301.6667
$ws.Range("M107").Value = -1381.6667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 18585.273
$ws.Range("I7").Value = 28812.572
$ws.Range("J7").Value = 687.5
$ws.Range("K7").Value = 28812.572
$ws.Range("L7").Value = 687.5
$ws.Range("M7").Value = -28699.572
$ws.Range("N7").Value = -913.5
$ws.Range("H62").Value = 6682.4
$ws.Range("I62").Value = 5890
$ws.Range("J62").Value = 8531.333000000001
$ws.Range("K62").Value = 5890
$ws.Range("L62").Value = 8531.333000000001
$ws.Range("M62").Value = -5266
$ws.Range("N62").Value = -9779.333000000001
$ws.Range("H65").Value = 6682.4
$ws.Range("I65").Value = 5890
$ws.Range("J65").Value = 8531.333000000001
$ws.Range("K65").Value = 29450
$ws.Range("L65").Value = 42656.665
$ws.Range("M65").Value = -26330
$ws.Range("N65").Value = -48896.665
$ws.Range("H134").Value = 6286.9287
$ws.Range("J134").Value = 2310.9092
$ws.Range("L134").Value = 6932.7276
$ws.Range("N134").Value = -12002.7276

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1887.32
$ws.Range("I131").Value = 1599.6666
$ws.Range("J131").Value = 1896.2166
$ws.Range("K131").Value = 4798.9998
$ws.Range("L131").Value = 5688.6498
$ws.Range("M131").Value = 241.0002000000004
$ws.Range("N131").Value = -15768.6498

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.61539
$ws.Range("I2").Value = 114.347824
$ws.Range("J2").Value = 298.66666
$ws.Range("K2").Value = 114.347824
$ws.Range("L2").Value = 298.66666
$ws.Range("M2").Value = -1.347824000000003
$ws.Range("N2").Value = -524.66666
$ws.Range("H70").Value = 8022
$ws.Range("I70").Value = 6782.636
$ws.Range("J70").Value = 10294.167
$ws.Range("K70").Value = 6782.636
$ws.Range("L70").Value = 10294.167
$ws.Range("M70").Value = -6512.636
$ws.Range("N70").Value = -10834.167
$ws.Range("H73").Value = 8022
$ws.Range("I73").Value = 6782.636
$ws.Range("J73").Value = 10294.167
$ws.Range("K73").Value = 6782.636
$ws.Range("L73").Value = 10294.167
$ws.Range("M73").Value = -5846.636
$ws.Range("N73").Value = -12166.167
$ws.Range("H80").Value = 10475.429
$ws.Range("J80").Value = 3911
$ws.Range("L80").Value = 3911
$ws.Range("N80").Value = -5907
$ws.Range("H83").Value = 10475.429
$ws.Range("J83").Value = 3911
$ws.Range("L83").Value = 19555
$ws.Range("N83").Value = -29539
$ws.Range("H97").Value = 11709.75
$ws.Range("I97").Value = 14178.556
$ws.Range("J97").Value = 4303.3335
$ws.Range("K97").Value = 14178.556
$ws.Range("L97").Value = 4303.3335
$ws.Range("M97").Value = -13682.556
$ws.Range("N97").Value = -5295.3335
$ws.Range("H102").Value = 9509.105
$ws.Range("I102").Value = 11539.615
$ws.Range("K102").Value = 11539.615
$ws.Range("M102").Value = -9917.615
$ws.Range("H122").Value = 15010.786
$ws.Range("J122").Value = 13888.333
$ws.Range("L122").Value = 41664.999
$ws.Range("N122").Value = -46564.999
$ws.Range("H126").Value = 8508.762000000001
$ws.Range("J126").Value = 3416.182
$ws.Range("L126").Value = 10248.546
$ws.Range("N126").Value = -15188.546

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25727.38
$ws.Range("J40").Value = 9323.5
$ws.Range("L40").Value = 9323.5
$ws.Range("N40").Value = -9595.5
$ws.Range("H93").Value = 5745.5884
$ws.Range("I93").Value = 7566.4165
$ws.Range("K93").Value = 7566.4165
$ws.Range("M93").Value = -6318.4165
$ws.Range("H100").Value = 16571.285
$ws.Range("I100").Value = 26333
$ws.Range("K100").Value = 26333
$ws.Range("M100").Value = -25792
$ws.Range("H122").Value = 5281.7144
$ws.Range("I122").Value = 5289.4736
$ws.Range("K122").Value = 15868.4208
$ws.Range("M122").Value = -13418.4208
$ws.Range("H136").Value = 3798.0312
$ws.Range("I136").Value = 2795.4
$ws.Range("J136").Value = 5469.0835
$ws.Range("K136").Value = 8386.200000000001
$ws.Range("L136").Value = 16407.2505
$ws.Range("M136").Value = -5836.200000000001
$ws.Range("N136").Value = -21507.2505

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
